$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it
# stored as text (matches the source data, which is all inline strings),
# then strip the temporary text number-format so no extra cell style lingers.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Rows 35/36: Hedera and Celestia swap list positions, with refreshed price/volume
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D35" "18.88"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D36" "0.0790"
$ws.Range("E36").Value = "  +0.05%  "

# Price (D) and Volume(1h) (E) refresh for the remaining rows
Set-TextValue "D2" "42.883.68"
$ws.Range("E2").Value = "  +0.32%  "
Set-TextValue "D3" "2.526.81"
$ws.Range("E3").Value = "  +0.18%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue "D5" "317.00"
$ws.Range("E5").Value = "  +4.01%  "
Set-TextValue "D6" "95.20"
$ws.Range("E6").Value = "  -1.56%  "
Set-TextValue "D7" "0.578"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.06%  "
Set-TextValue "D9" "0.532"
$ws.Range("E9").Value = "  -1.44%  "
Set-TextValue "D10" "36.00"
$ws.Range("E10").Value = "  -1.15%  "
Set-TextValue "D11" "0.0810"
$ws.Range("E11").Value = "  -0.08%  "
Set-TextValue "D12" "7.57"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("E13").Value = "  -0.54%  "
Set-TextValue "D14" "2.916.16"
$ws.Range("E14").Value = "  +0.28%  "
Set-TextValue "D15" "2.527.77"
$ws.Range("E15").Value = "  +0.76%  "
Set-TextValue "D16" "15.26"
$ws.Range("E16").Value = "  +1.13%  "
Set-TextValue "D17" "0.848"
$ws.Range("E17").Value = "  -1.38%  "
Set-TextValue "D18" "42.956.79"
$ws.Range("E18").Value = "  +0.56%  "
Set-TextValue "D19" "12.95"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  +3.06%  "
Set-TextValue "D21" "0.0₃0965"
$ws.Range("E21").Value = "  -0.99%  "
Set-TextValue "D22" "70.11"
$ws.Range("E22").Value = "  -1.55%  "
Set-TextValue "D23" "251.37"
$ws.Range("E23").Value = "  +0.01%  "
Set-TextValue "D24" "2.95"
$ws.Range("E24").Value = "  +1.25%  "
Set-TextValue "D25" "2.01"
$ws.Range("E25").Value = "  -1.28%  "
Set-TextValue "D26" "26.87"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +3.64%  "
Set-TextValue "D29" "39.86"
$ws.Range("E29").Value = "  +4.19%  "
Set-TextValue "D30" "10.26"
$ws.Range("E30").Value = "  -0.51%  "
Set-TextValue "D31" "6.02"
$ws.Range("E31").Value = "  +0.84%  "
Set-TextValue "D32" "154.64"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("E38").Value = "  -3.44%  "
$ws.Range("E39").Value = "  -0.36%  "
Set-TextValue "D40" "23.74"
$ws.Range("E40").Value = "  -1.69%  "
Set-TextValue "D41" "2.22"
$ws.Range("E41").Value = "  +7.90%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  +0.27%  "
Set-TextValue "D44" "3.77"
$ws.Range("E44").Value = "  -2.04%  "
Set-TextValue "D45" "3.28"
$ws.Range("E45").Value = "  -3.37%  "
Set-TextValue "D46" "2.018.70"
$ws.Range("E46").Value = "  -0.73%  "
Set-TextValue "D47" "85.86"
$ws.Range("E47").Value = "  +0.82%  "
Set-TextValue "D48" "8.78"
$ws.Range("E48").Value = "  -1.56%  "
Set-TextValue "D49" "2.771.16"
$ws.Range("E49").Value = "  +0.13%  "
Set-TextValue "D50" "73.38"
$ws.Range("E50").Value = "  +1.79%  "
Set-TextValue "D51" "102.40"
$ws.Range("E51").Value = "  +0.55%  "
